$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "http://testing.bsbtest.com/default.aspx?portalid=1369 "
$ws.Range("B3").Value = "host"
$ws.Range("B4").Value = "fnw00t#"

$ws.Range("D18").Select()
